$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-mark numeric-looking text cells as Text format so their values are not
# auto-converted to numbers when assigned below.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '27.914.39'
$ws.Range('E2').Value = '  -0.32%  '

$ws.Range('D3').Value = '1.910.15'
$ws.Range('E3').Value = '  +0.18%  '

$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').Value = '  -0.51%  '

$ws.Range('D5').Value = '312.82'
$ws.Range('E5').Value = '  -1.56%  '

$ws.Range('D6').Value = '0.9991'

$ws.Range('D7').Value = '0.4993'
$ws.Range('E7').Value = '  +3.47%  '

$ws.Range('D8').Value = '0.3807'
$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').Value = '0.07285'
$ws.Range('E9').Value = '  -1.05%  '

$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value = '21.32'
$ws.Range('E10').Value = '  +2.50%  '

$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value = '0.9091'
$ws.Range('E11').Value = '  -2.56%  '

$ws.Range('D12').Value = '0.07643'
$ws.Range('E12').Value = '  -1.62%  '

$ws.Range('D13').Value = '1.870.69'
$ws.Range('E13').Value = '  -1.80%  '

$ws.Range('D14').Value = '5.475'
$ws.Range('E14').Value = '  -0.42%  '

$ws.Range('D15').Value = '92.62'
$ws.Range('E15').Value = '  +0.86%  '

$ws.Range('D16').Value = '1.0000'
$ws.Range('E16').Value = '  -0.60%  '

$ws.Range('D17').Value = '0.000008724'
$ws.Range('E17').Value = '  -1.88%  '

$ws.Range('D18').Value = '0.9985'
$ws.Range('E18').Value = '  -0.58%  '

$ws.Range('D19').Value = '27.943.78'
$ws.Range('E19').Value = '  -0.32%  '

$ws.Range('D20').Value = '14.65'
$ws.Range('E20').Value = '  -0.61%  '

$ws.Range('D21').Value = '5.166'
$ws.Range('E21').Value = '  +0.21%  '

$ws.Range('D22').Value = '2.110.88'
$ws.Range('E22').Value = '  -1.25%  '

$ws.Range('D23').Value = '10.87'
$ws.Range('E23').Value = '  -0.30%  '

$ws.Range('D24').Value = '6.605'
$ws.Range('E24').Value = '  -0.43%  '

$ws.Range('D25').Value = '152.99'
$ws.Range('E25').Value = '  -2.57%  '

$ws.Range('D26').Value = '1.840'
$ws.Range('E26').Value = '  -3.70%  '

$ws.Range('D27').Value = '2.221'
$ws.Range('E27').Value = '  +5.02%  '

$ws.Range('D28').Value = '18.39'
$ws.Range('E28').Value = '  -0.62%  '

$ws.Range('D29').Value = '115.05'
$ws.Range('E29').Value = '  -1.94%  '

$ws.Range('D30').Value = '4.903'
$ws.Range('E30').Value = '  -1.67%  '

$ws.Range('D31').Value = '0.09002'
$ws.Range('E31').Value = '  +0.65%  '

$ws.Range('E32').Value = '  -1.89%  '

$ws.Range('D33').Value = '4.822'
$ws.Range('E33').Value = '  +3.41%  '

$ws.Range('E34').Value = '  -1.64%  '

$ws.Range('D35').Value = '0.7795'
$ws.Range('E35').Value = '  +1.02%  '

$ws.Range('D36').Value = '2.629'
$ws.Range('E36').Value = '  +1.31%  '

$ws.Range('D37').Value = '0.02082'
$ws.Range('E37').Value = '  +1.37%  '

$ws.Range('D38').Value = '3.058'
$ws.Range('E38').Value = '  +2.14%  '

$ws.Range('D39').Value = '1.091'
$ws.Range('E39').Value = '  -1.56%  '

$ws.Range('D40').Value = '0.5549'
$ws.Range('E40').Value = '  +0.74%  '

$ws.Range('D41').Value = '0.05271'
$ws.Range('E41').Value = '  -0.44%  '

$ws.Range('D42').Value = '6.816'
$ws.Range('E42').Value = '  -2.49%  '

$ws.Range('D43').Value = '114.07'
$ws.Range('E43').Value = '  +3.73%  '

$ws.Range('D44').Value = '8.504'
$ws.Range('E44').Value = '  +0.22%  '

$ws.Range('E45').Value = '  -0.70%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '10.62'
$ws.Range('E46').Value = '  -0.73%  '

$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.4820'
$ws.Range('E47').Value = '  -0.03%  '

$ws.Range('D48').Value = '0.9990'
$ws.Range('E48').Value = '  -0.56%  '

$ws.Range('D49').Value = '1.640'
$ws.Range('E49').Value = '  -0.42%  '

$ws.Range('D50').Value = '67.24'
$ws.Range('E50').Value = '  -1.18%  '

$ws.Range('D51').Value = '0.06049'
$ws.Range('E51').Value = '  -0.44%  '

# Restore original (default) cell formatting on the cells we temporarily
# switched to Text format, so the saved style table matches the original.
$ws.Range('D4').ClearFormats()
$ws.Range('D5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('D7').ClearFormats()
$ws.Range('D8').ClearFormats()
$ws.Range('D9').ClearFormats()
$ws.Range('D10').ClearFormats()
$ws.Range('D11').ClearFormats()
$ws.Range('D12').ClearFormats()
$ws.Range('D14').ClearFormats()
$ws.Range('D15').ClearFormats()
$ws.Range('D16').ClearFormats()
$ws.Range('D17').ClearFormats()
$ws.Range('D18').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D23').ClearFormats()
$ws.Range('D24').ClearFormats()
$ws.Range('D25').ClearFormats()
$ws.Range('D26').ClearFormats()
$ws.Range('D27').ClearFormats()
$ws.Range('D28').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('D30').ClearFormats()
$ws.Range('D31').ClearFormats()
$ws.Range('D33').ClearFormats()
$ws.Range('D35').ClearFormats()
$ws.Range('D36').ClearFormats()
$ws.Range('D37').ClearFormats()
$ws.Range('D38').ClearFormats()
$ws.Range('D39').ClearFormats()
$ws.Range('D40').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('D42').ClearFormats()
$ws.Range('D43').ClearFormats()
$ws.Range('D44').ClearFormats()
$ws.Range('D46').ClearFormats()
$ws.Range('D47').ClearFormats()
$ws.Range('D48').ClearFormats()
$ws.Range('D49').ClearFormats()
$ws.Range('D50').ClearFormats()
$ws.Range('D51').ClearFormats()
